# Update cryptos list: refresh Price / Volume(1h) figures and re-rank a few coins
# that changed relative order, per the upstream GitHub Actions data pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Sheet, [string]$Addr, [string]$Val)
    $cell = $Sheet.Range($Addr)
    # Force text storage so numeric-looking strings (e.g. '139.10', '4.35')
    # are not silently coerced into floating point numbers, and so that
    # thousand-dot-separated price strings (e.g. '55.176.47') stay verbatim.
    $cell.NumberFormat = '@'
    $cell.Value = $Val
    # Drop back to the default style so no stray number-format style sticks
    # around on the cell (matches the original unstyled D/E cells).
    $cell.Style = 'Normal'
}

Set-TextCell $ws 'D2' '55.176.47'
Set-TextCell $ws 'E2' '  +6.31%  '
Set-TextCell $ws 'D3' '2.435.23'
Set-TextCell $ws 'E3' '  +6.72%  '
Set-TextCell $ws 'E4' '  -0.12%  '
Set-TextCell $ws 'D5' '480.61'
Set-TextCell $ws 'E5' '  +9.43%  '
Set-TextCell $ws 'D6' '139.10'
Set-TextCell $ws 'E6' '  +17.42%  '
Set-TextCell $ws 'D7' '0.997'
Set-TextCell $ws 'E7' '  -0.15%  '
Set-TextCell $ws 'D8' '0.502'
Set-TextCell $ws 'E8' '  +8.93%  '
Set-TextCell $ws 'D9' '2.458.09'
Set-TextCell $ws 'E9' '  +7.84%  '
Set-TextCell $ws 'D10' '0.0968'
Set-TextCell $ws 'E10' '  +12.82%  '
Set-TextCell $ws 'E11' '  +4.54%  '
Set-TextCell $ws 'D12' '0.325'
Set-TextCell $ws 'E12' '  +9.40%  '
Set-TextCell $ws 'E13' '  +2.38%  '
Set-TextCell $ws 'D14' '2.869.14'
Set-TextCell $ws 'E14' '  +7.37%  '
Set-TextCell $ws 'D15' '55.096.97'
Set-TextCell $ws 'E15' '  +6.05%  '
Set-TextCell $ws 'D16' '20.48'
Set-TextCell $ws 'E16' '  +10.40%  '
Set-TextCell $ws 'D17' '0.0000135'
Set-TextCell $ws 'E17' '  +16.65%  '
Set-TextCell $ws 'D18' '2.453.76'
Set-TextCell $ws 'E18' '  +6.38%  '
Set-TextCell $ws 'D19' '4.35'
Set-TextCell $ws 'E19' '  +11.01%  '
Set-TextCell $ws 'B20' 'Chainlink'
Set-TextCell $ws 'C20' 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextCell $ws 'D20' '9.93'
Set-TextCell $ws 'E20' '  +14.58%  '
Set-TextCell $ws 'B21' 'BitcoinCash'
Set-TextCell $ws 'C21' 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextCell $ws 'D21' '315.11'
Set-TextCell $ws 'E21' '  +7.21%  '
Set-TextCell $ws 'D22' '0.996'
Set-TextCell $ws 'E22' '  +0.46%  '
Set-TextCell $ws 'D23' '5.63'
Set-TextCell $ws 'E23' '  +11.41%  '
Set-TextCell $ws 'D24' '57.22'
Set-TextCell $ws 'E24' '  +8.33%  '
Set-TextCell $ws 'E25' '  +11.51%  '
Set-TextCell $ws 'E26' '  +0.19%  '
Set-TextCell $ws 'E27' '  +16.46%  '
Set-TextCell $ws 'D28' '2.547.80'
Set-TextCell $ws 'E28' '  +6.62%  '
Set-TextCell $ws 'D29' '7.34'
Set-TextCell $ws 'E29' '  +7.63%  '
Set-TextCell $ws 'E30' '  +21.67%  '
Set-TextCell $ws 'D31' '0.998'
Set-TextCell $ws 'E31' '  -0.05%  '
Set-TextCell $ws 'D32' '148.67'
Set-TextCell $ws 'E32' '  +4.28%  '
Set-TextCell $ws 'D33' '17.96'
Set-TextCell $ws 'E33' '  +8.24%  '
Set-TextCell $ws 'E34' '  +12.56%  '
Set-TextCell $ws 'D35' '5.15'
Set-TextCell $ws 'E35' '  +11.41%  '
Set-TextCell $ws 'E36' '  +14.57%  '
Set-TextCell $ws 'B37' 'Fetch.AI'
Set-TextCell $ws 'C37' 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextCell $ws 'D37' '0.857'
Set-TextCell $ws 'E37' '  +7.31%  '
Set-TextCell $ws 'B38' 'NEARProtocol'
Set-TextCell $ws 'C38' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextCell $ws 'D38' '3.60'
Set-TextCell $ws 'E38' '  +7.98%  '
Set-TextCell $ws 'B39' 'FirstDigitalUSD'
Set-TextCell $ws 'C39' 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextCell $ws 'D39' '0.995'
Set-TextCell $ws 'E39' '  +0.00%  '
Set-TextCell $ws 'B40' 'OKB'
Set-TextCell $ws 'C40' 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextCell $ws 'D40' '33.19'
Set-TextCell $ws 'E40' '  +4.72%  '
Set-TextCell $ws 'B41' 'Mantle'
Set-TextCell $ws 'C41' 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextCell $ws 'D41' '0.602'
Set-TextCell $ws 'E41' '  +9.76%  '
Set-TextCell $ws 'B42' 'Filecoin'
Set-TextCell $ws 'C42' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextCell $ws 'D42' '3.40'
Set-TextCell $ws 'E42' '  +10.80%  '
Set-TextCell $ws 'D43' '0.0543'
Set-TextCell $ws 'E43' '  +9.40%  '
Set-TextCell $ws 'E44' '  +12.57%  '
Set-TextCell $ws 'D45' '4.65'
Set-TextCell $ws 'E45' '  +14.72%  '
Set-TextCell $ws 'D46' '10.10'
Set-TextCell $ws 'E46' '  -0.13%  '
Set-TextCell $ws 'D47' '254.41'
Set-TextCell $ws 'E47' '  +30.40%  '
Set-TextCell $ws 'D48' '0.0898'
Set-TextCell $ws 'E48' '  +10.83%  '
Set-TextCell $ws 'B49' 'VeChain'
Set-TextCell $ws 'C49' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextCell $ws 'D49' '0.0222'
Set-TextCell $ws 'E49' '  +10.42%  '
Set-TextCell $ws 'B50' 'Maker'
Set-TextCell $ws 'C50' 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextCell $ws 'D50' '1.931.07'
Set-TextCell $ws 'E50' '  +1.83%  '
Set-TextCell $ws 'D51' '17.08'
Set-TextCell $ws 'E51' '  +10.75%  '
